$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 106.14286
$ws.Range("I33").Value = 47.75
$ws.Range("J33").Value = 184
$ws.Range("K33").Value = 47.75
$ws.Range("L33").Value = 184
$ws.Range("M33").Value = 181.25
$ws.Range("N33").Value = -642
$ws.Range("H76").Value = 111115640
$ws.Range("I76").Value = 2186.75
$ws.Range("J76").Value = 200006400
$ws.Range("K76").Value = 2186.75
$ws.Range("L76").Value = 200006400
$ws.Range("M76").Value = -1871.75
$ws.Range("N76").Value = -200007030
$ws.Range("H79").Value = 111115640
$ws.Range("I79").Value = 2186.75
$ws.Range("J79").Value = 200006400
$ws.Range("K79").Value = 2186.75
$ws.Range("L79").Value = 200006400
$ws.Range("M79").Value = -1094.75
$ws.Range("N79").Value = -200008584
$ws.Range("H98").Value = 1545.2222
$ws.Range("I98").Value = 934.5
$ws.Range("K98").Value = 934.5
$ws.Range("M98").Value = 563.5
$ws.Range("H122").Value = 1545.2222
$ws.Range("I122").Value = 934.5
$ws.Range("K122").Value = 2803.5
$ws.Range("M122").Value = -353.5
$ws.Range("H132").Value = 3942.7273
$ws.Range("I132").Value = 4738.2856
$ws.Range("J132").Value = 2550.5
$ws.Range("K132").Value = 14214.8568
$ws.Range("L132").Value = 7651.5
$ws.Range("M132").Value = -11684.8568
$ws.Range("N132").Value = -12711.5
$ws.Range("H137").Value = 1006965
$ws.Range("J137").Value = 1122761.5
$ws.Range("L137").Value = 3368284.5
$ws.Range("N137").Value = -3373384.5
$ws.Range("H138").Value = 7964.271
$ws.Range("I138").Value = 4225.25
$ws.Range("K138").Value = 12675.75
$ws.Range("M138").Value = -7535.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2911.0715
$ws.Range("I45").Value = 2405
$ws.Range("K45").Value = 2405
$ws.Range("M45").Value = -2028
$ws.Range("H63").Value = 5597.909
$ws.Range("I63").Value = 3514.5
$ws.Range("K63").Value = 3514.5
$ws.Range("M63").Value = -2828.5
$ws.Range("H66").Value = 5597.909
$ws.Range("I66").Value = 3514.5
$ws.Range("K66").Value = 17572.5
$ws.Range("M66").Value = -14140.5
$ws.Range("H74").Value = 2710.875
$ws.Range("I74").Value = 2710.875
$ws.Range("K74").Value = 2710.875
$ws.Range("M74").Value = -1836.875
$ws.Range("H77").Value = 2710.875
$ws.Range("I77").Value = 2710.875
$ws.Range("K77").Value = 13554.375
$ws.Range("M77").Value = -9186.375
$ws.Range("H132").Value = 3994.3635
$ws.Range("I132").Value = 2988.3235
$ws.Range("K132").Value = 8964.970499999999
$ws.Range("M132").Value = -6434.970499999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 55984.316
$ws.Range("I134").Value = 2934.923
$ws.Range("K134").Value = 8804.769
$ws.Range("M134").Value = -6269.769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36246.78
$ws.Range("I31").Value = 1722.8667
$ws.Range("J31").Value = 66709.06
$ws.Range("K31").Value = 1722.8667
$ws.Range("L31").Value = 66709.06
$ws.Range("M31").Value = -1427.8667
$ws.Range("N31").Value = -67299.06
$ws.Range("H34").Value = 36246.78
$ws.Range("I34").Value = 1722.8667
$ws.Range("J34").Value = 66709.06
$ws.Range("K34").Value = 1722.8667
$ws.Range("L34").Value = 66709.06
$ws.Range("M34").Value = -1520.8667
$ws.Range("N34").Value = -67113.06
$ws.Range("H58").Value = 506616.5
$ws.Range("I58").Value = 838869.2
$ws.Range("J58").Value = 8237.5
$ws.Range("K58").Value = 838869.2
$ws.Range("L58").Value = 8237.5
$ws.Range("M58").Value = -838666.2
$ws.Range("N58").Value = -8643.5
$ws.Range("H86").Value = 83336
$ws.Range("I86").Value = 50000
$ws.Range("K86").Value = 50000
$ws.Range("M86").Value = -48877
$ws.Range("H89").Value = 83336
$ws.Range("I89").Value = 50000
$ws.Range("K89").Value = 250000
$ws.Range("M89").Value = -244384
$ws.Range("H107").Value = 242.3
$ws.Range("I107").Value = 242.3
$ws.Range("K107").Value = 242.3
$ws.Range("M107").Value = 1677.7
$ws.Range("H122").Value = 4202.3125
$ws.Range("I122").Value = 2943.111
$ws.Range("J122").Value = 5821.2856
$ws.Range("K122").Value = 8829.332999999999
$ws.Range("L122").Value = 17463.8568
$ws.Range("M122").Value = -6379.332999999999
$ws.Range("N122").Value = -22363.8568
$ws.Range("H132").Value = 4948.846
$ws.Range("I132").Value = 4610.8945
$ws.Range("J132").Value = 5866.143
$ws.Range("K132").Value = 13832.6835
$ws.Range("L132").Value = 17598.429
$ws.Range("M132").Value = -11302.6835
$ws.Range("N132").Value = -22658.429
$ws.Range("H134").Value = 2227367.5
$ws.Range("I134").Value = 1434115.4
$ws.Range("K134").Value = 4302346.199999999
$ws.Range("M134").Value = -4299811.199999999
$ws.Range("H136").Value = 506616.5
$ws.Range("I136").Value = 838869.2
$ws.Range("J136").Value = 8237.5
$ws.Range("K136").Value = 2516607.6
$ws.Range("L136").Value = 24712.5
$ws.Range("M136").Value = -2514057.6
$ws.Range("N136").Value = -29812.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1200213.1
$ws.Range("I5").Value = 100239.125
$ws.Range("K5").Value = 300717.375
$ws.Range("M5").Value = -300605.375
$ws.Range("H63").Value = 1628
$ws.Range("I63").Value = 1837.3334
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 5512.0002
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -4763.0002
$ws.Range("N63").Value = -4498
$ws.Range("H66").Value = 1628
$ws.Range("I66").Value = 1837.3334
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 16536.0006
$ws.Range("L66").Value = 9000
$ws.Range("M66").Value = -12792.0006
$ws.Range("N66").Value = -16488
$ws.Range("H69").Value = 2625
$ws.Range("J69").Value = 2625
$ws.Range("L69").Value = 7875
$ws.Range("N69").Value = -9497
$ws.Range("H72").Value = 2625
$ws.Range("J72").Value = 2625
$ws.Range("L72").Value = 23625
$ws.Range("N72").Value = -31737
$ws.Range("H86").Value = 499
$ws.Range("J86").Value = 499
$ws.Range("L86").Value = 1497
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 499
$ws.Range("J89").Value = 499
$ws.Range("L89").Value = 4491
$ws.Range("N89").Value = -16347
$ws.Range("H122").Value = 56311.332
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 112122.664
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 1009103.976
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -1014003.976
$ws.Range("H129").Value = 1746.6
$ws.Range("I129").Value = 700
$ws.Range("K129").Value = 2100
$ws.Range("M129").Value = 2900
$ws.Range("H131").Value = 12423460
$ws.Range("I131").Value = 25718592
$ws.Range("J131").Value = 77981.36
$ws.Range("K131").Value = 77155776
$ws.Range("L131").Value = 233944.08
$ws.Range("M131").Value = -77150736
$ws.Range("N131").Value = -244024.08
$ws.Range("H132").Value = 442806.97
$ws.Range("I132").Value = 92292.45
$ws.Range("J132").Value = 718211.2
$ws.Range("K132").Value = 830632.0499999999
$ws.Range("L132").Value = 6463900.8
$ws.Range("M132").Value = -828102.0499999999
$ws.Range("N132").Value = -6468960.8
$ws.Range("H135").Value = 1200213.1
$ws.Range("I135").Value = 100239.125
$ws.Range("K135").Value = 902152.125
$ws.Range("M135").Value = -899617.125
$ws.Range("H137").Value = 2059.2727
$ws.Range("I137").Value = 1183.5555
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 3550.6665
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = 1549.3335
$ws.Range("N137").Value = -28200
$ws.Range("H139").Value = 6503.6
$ws.Range("I139").Value = 6503.6
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 19510.8
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -14370.8
$ws.Range("N139").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 59500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 59500
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -60518
$ws.Range("H126").Value = 4165.8887
$ws.Range("J126").Value = 4889.6665
$ws.Range("L126").Value = 14668.9995
$ws.Range("N126").Value = -19608.9995
$ws.Range("H132").Value = 920610
$ws.Range("I132").Value = 1254027.5
$ws.Range("J132").Value = 253775
$ws.Range("K132").Value = 3762082.5
$ws.Range("L132").Value = 761325
$ws.Range("M132").Value = -3759552.5
$ws.Range("N132").Value = -766385

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4999
$ws.Range("I2").Value = 4999
$ws.Range("K2").Value = 4999
$ws.Range("M2").Value = -4887
$ws.Range("H46").Value = 4828.75
$ws.Range("I46").Value = 4968.1875
$ws.Range("J46").Value = 4549.875
$ws.Range("K46").Value = 4968.1875
$ws.Range("L46").Value = 4549.875
$ws.Range("M46").Value = -4780.1875
$ws.Range("N46").Value = -4925.875
$ws.Range("H132").Value = 3968.75
$ws.Range("I132").Value = 1412.5
$ws.Range("J132").Value = 6525
$ws.Range("K132").Value = 4237.5
$ws.Range("L132").Value = 19575
$ws.Range("M132").Value = -1707.5
$ws.Range("N132").Value = -24635
$ws.Range("H136").Value = 1256782.2
$ws.Range("I136").Value = 2227668
$ws.Range("J136").Value = 8500.429
$ws.Range("K136").Value = 6683004
$ws.Range("L136").Value = 25501.287
$ws.Range("M136").Value = -6680454
$ws.Range("N136").Value = -30601.287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 39827.715
$ws.Range("I132").Value = 4008.8235
$ws.Range("K132").Value = 12026.4705
$ws.Range("M132").Value = -9496.470499999999
